$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-26 Saturday", "2025-07-27 Sunday"),
    @("97×73=7081", "38×11=418"),
    @("72×13=936", "62×16=992"),
    @("42×72=3024", "12×97=1164"),
    @("88×30=2640", "55×28=1540"),
    @("15×89=1335", "63×29=1827"),
    @("76×15=1140", "44×83=3652"),
    @("22×71=1562", "21×14=294"),
    @("11×18=198", "57×96=5472"),
    @("83×16=1328", "23×41=943"),
    @("93×25=2325", "99×22=2178"),
    @("81×60=4860", "68×72=4896"),
    @("66×61=4026", "56×85=4760"),
    @("32×56=1792", "50×41=2050"),
    @("66×88=5808", "76×41=3116"),
    @("60×91=5460", "40×82=3280"),
    @("74×64=4736", "58×67=3886"),
    @("66×71=4686", "99×94=9306"),
    @("77×26=2002", "11×69=759"),
    @("49×52=2548", "41×95=3895"),
    @("60×45=2700", "77×16=1232"),
    @("46×46=2116", "50×58=2900"),
    @("58×64=3712", "52×74=3848"),
    @("13×71=923", "56×91=5096"),
    @("93×88=8184", "75×45=3375"),
    @("96×24=2304", "28×75=2100")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
